$d = $word.ActiveDocument

function Split-IntoRuns($startPos, $segments) {
    # Toggle a formatting property on/off on each sub-range in turn; this
    # forces Word to break the run at each boundary without altering the
    # visible formatting of the final text.
    $pos = $startPos
    foreach ($seg in $segments) {
        $segLen = $seg.Length
        if ($segLen -gt 0) {
            $r = $d.Range($pos, $pos + $segLen)
            $r.Bold = 1
            $r.Bold = 0
        }
        $pos = $pos + $segLen
    }
}

# ---------------------------------------------------------------------
# 1) Phone number: "+1-8572609294" -> "+1-(857)-260-9294" split across
#    several runs.
# ---------------------------------------------------------------------
$oldPhone = "+1-8572609294"
$newPhone = "+1-(857)-260-9294"
$rng = $d.Content
$found = $rng.Find.Execute($oldPhone)
if ($found) {
    $start = $rng.Start
    $rng.Text = $newPhone
    $phoneSegs = @("+1-", "(", "857", ")-", "260", "-", "9294")
    Split-IntoRuns $start $phoneSegs
}

# ---------------------------------------------------------------------
# 2) "Predicted Next Best Action ..." bullet rewritten and split across
#    several runs.
# ---------------------------------------------------------------------
$oldPnba = "Predicted Next Best Action for a coupon generator application using K-Means Clustering with an accuracy of 61 %."
$newPnba = "Predicted Next Best Action for an Offer Generator using K-Means Clustering with a 61 % average chance of achieving the intents."
$rng = $d.Content
$found = $rng.Find.Execute($oldPnba)
if ($found) {
    $start = $rng.Start
    $rng.Text = $newPnba
    $pnbaSegs = @(
        "Predicted Next Best Action for a",
        "n",
        " ",
        "Offer",
        " ",
        "G",
        "enerator using K-Means Clustering with a 61 %",
        " average chance of achieving the intents",
        "."
    )
    Split-IntoRuns $start $pnbaSegs
}

# ---------------------------------------------------------------------
# 3) "Automated resume matching process ..." bullet: two runs merged
#    into a single run with updated wording.
# ---------------------------------------------------------------------
$oldResume = "Automated resume matching process using an NLP model and decreased the time spent by recruiting by approximately 80 %."
$newResume = "Automated resume matching process using a word count model and decreased the time spent by recruiting by ~ 80 %."
$rng = $d.Content
$rng.Find.Execute($oldResume, $true, $false, $false, $false, $false, $true, 1, $false, $newResume, 2)

Write-Output "done"
